# NetConnect-NetChain2 TestData.xlsx -- "added values for Vendor"
#
# Row 9  (NetchainTest.CreateVendor)      : fill in the remaining Vendor fields (G:AB)
# Row 10 (NEW - NetchainTest.CreateVendorNeg): duplicate vendor row used for the negative test,
#                                              inserted above the old row 10, pushing
#                                              "NetchainTest.CreateGoodsAndServices" down to row 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at 10 so the old row 10 (CreateGoodsAndServices) becomes row 11.
# ---------------------------------------------------------------------------
$ws.Range("A10").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2. Finish populating row 9 (NetchainTest.CreateVendor) -- columns A:F already
#    hold the right values, only G:AB are new/changed.
# ---------------------------------------------------------------------------
$ws.Range("G9").Value = "'9098989988"
$ws.Range("H9").Value = "Tester"
$ws.Range("I9").Value = "pune"
$ws.Range("J9").Value = "Department1"
$ws.Range("K9").Value = "Fuel"
$ws.Range("L9").Value = "Ms"
$ws.Range("M9").Value = "Abc"
$ws.Range("N9").Value = "Xyz"
$ws.Range("O9").Value = "Company1"
$ws.Range("P9").Value = "Company1"
$ws.Range("Q9").Value = "ltd"
$ws.Range("R9").Value = "other"
$ws.Range("S9").Value = "whc road"
$ws.Range("T9").Value = "'nagpur"
$ws.Range("T9").NumberFormat = "0"
$ws.Range("U9").Value = "'mh"
$ws.Range("U9").NumberFormat = "0"
$ws.Range("V9").Value = 440015
$ws.Range("W9").Value = "abcd@gmail.com"
$ws.Range("X9").Value = "'9098987766"
$ws.Range("Y9").Value = "'8989887677"
$ws.Range("Z9").Value = 12345
$ws.Range("AA9").Value = "www.abcd.com"
$ws.Range("AB9").Value = "notes"

# ---------------------------------------------------------------------------
# 3. Populate the new row 10 (NetchainTest.CreateVendorNeg) -- same vendor
#    payload as row 9, re-used for a negative-path test; W10 is entered as a
#    plain number (2) instead of the email address.
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "NetchainTest.CreateVendorNeg"
$ws.Range("B10").Value = "company1"
$ws.Range("C10").Value = "mission1"
$ws.Range("D10").Value = "mission&product1"
$ws.Range("E10").Value = "Vendor 1"
$ws.Range("F10").Value = "abcd@gmail.com"
$ws.Range("G10").Value = "'9098989988"
$ws.Range("H10").Value = "Tester"
$ws.Range("I10").Value = "pune"
$ws.Range("J10").Value = "Department1"
$ws.Range("K10").Value = "Fuel"
$ws.Range("L10").Value = "Ms"
$ws.Range("M10").Value = "Abc"
$ws.Range("N10").Value = "Xyz"
$ws.Range("O10").Value = "Company1"
$ws.Range("P10").Value = "Company1"
$ws.Range("Q10").Value = "ltd"
$ws.Range("R10").Value = "other"
$ws.Range("S10").Value = "whc road"
$ws.Range("T10").Value = "'nagpur"
$ws.Range("T10").NumberFormat = "0"
$ws.Range("U10").Value = "'mh"
$ws.Range("U10").NumberFormat = "0"
$ws.Range("V10").Value = 440015
$ws.Range("W10").Value = 2
$ws.Range("X10").Value = "'9098987766"
$ws.Range("Y10").Value = "'8989887677"
$ws.Range("Z10").Value = 12345
$ws.Range("AA10").Value = "www.abcd.com"
$ws.Range("AB10").Value = "notes"

# Match the recorded row height / selection state for the edited rows.
$ws.Rows.Item(9).RowHeight = 35.25
$ws.Rows.Item(10).RowHeight = 35.25
$ws.Rows.Item(11).RowHeight = 35.25

$ws.Range("W10").Select()
